{"js": "// no-op to test baseline\n", "ps1": "# no-op to test baseline\n"}
